$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.721.98"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.778.34"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "437.23"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.14"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  -8.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000319"
$ws.Range("E11").Value = "  -12.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.81"
$ws.Range("E12").Value = "  +4.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.42"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "4.376.69"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.82"
$ws.Range("E15").Value = "  -5.19%  "
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "3.756.00"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.92"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("E19").Value = "  +5.76%  "
$ws.Range("D20").Value = "66.738.95"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "412.38"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.53"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  +7.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.10"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.13"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.67"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  +31.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "730.69"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  +9.60%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.76"
$ws.Range("E32").Value = "  +10.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.44"
$ws.Range("E34").Value = "  +12.18%  "
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("E37").Value = "  +24.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.53"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0479"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.64"
$ws.Range("E40").Value = "  +32.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0681"
$ws.Range("E42").Value = "  -16.08%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  +6.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.328"
$ws.Range("E46").Value = "  +12.53%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.66"
$ws.Range("E48").Value = "  +4.65%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.69"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("E51").Value = "  +1.46%  "
